$wb = $excel.ActiveWorkbook

# Rename Sheet3 -> addValidEmployeeTest
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Name = "addValidEmployeeTest"

# Populate the new employee test data
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "First Name"
$ws.Range("D1").Value = "Middle Name"
$ws.Range("E1").Value = "Last Name"

$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("C2").Value = "John"
$ws.Range("D2").Value = "w"
$ws.Range("E2").Value = "Wick"

$ws.Range("A3").Value = "Admin"
$ws.Range("B3").Value = "admin123"
$ws.Range("C3").Value = "Jack"
$ws.Range("D3").Value = "w"
$ws.Range("E3").Value = "Wick"

# Match the font used for the new data (Arial 10, explicit black)
$ws.Range("A1:E3").Font.Name = "Arial"
$ws.Range("A1:E3").Font.Size = 10
$ws.Range("A1:E3").Font.Color = 0

# Match the row heights used when the data was typed in
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75

$ws.Range("E3").Select()

# Page setup for the new sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Make the new sheet the active tab (moves tabSelected / activeTab too)
$ws.Activate()
